$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change A3 from "99" to "9".
# Leading apostrophe forces Excel to store it as text (quote-prefixed),
# matching the original cell's text formatting/style and adding a new
# shared-string entry "9" instead of re-using the "99" entry.
$ws.Range("A3").Value = "'9"
